$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).NumberFormat = '@'
$ws.Cells.Item(2,4).Value = '63.205.60'
$ws.Cells.Item(2,4).ClearFormats()
$ws.Cells.Item(2,5).NumberFormat = '@'
$ws.Cells.Item(2,5).Value = '  +1.28%  '
$ws.Cells.Item(2,5).ClearFormats()

# Row 3
$ws.Cells.Item(3,4).NumberFormat = '@'
$ws.Cells.Item(3,4).Value = '2.477.85'
$ws.Cells.Item(3,4).ClearFormats()
$ws.Cells.Item(3,5).NumberFormat = '@'
$ws.Cells.Item(3,5).Value = '  +3.09%  '
$ws.Cells.Item(3,5).ClearFormats()

# Row 4
$ws.Cells.Item(4,5).NumberFormat = '@'
$ws.Cells.Item(4,5).Value = '  -0.66%  '
$ws.Cells.Item(4,5).ClearFormats()

# Row 5
$ws.Cells.Item(5,4).NumberFormat = '@'
$ws.Cells.Item(5,4).Value = '576.74'
$ws.Cells.Item(5,4).ClearFormats()
$ws.Cells.Item(5,5).NumberFormat = '@'
$ws.Cells.Item(5,5).Value = '  +0.64%  '
$ws.Cells.Item(5,5).ClearFormats()

# Row 6
$ws.Cells.Item(6,4).NumberFormat = '@'
$ws.Cells.Item(6,4).Value = '146.62'
$ws.Cells.Item(6,4).ClearFormats()
$ws.Cells.Item(6,5).NumberFormat = '@'
$ws.Cells.Item(6,5).Value = '  +0.88%  '
$ws.Cells.Item(6,5).ClearFormats()

# Row 7
$ws.Cells.Item(7,5).NumberFormat = '@'
$ws.Cells.Item(7,5).Value = '  +0.14%  '
$ws.Cells.Item(7,5).ClearFormats()

# Row 8
$ws.Cells.Item(8,4).NumberFormat = '@'
$ws.Cells.Item(8,4).Value = '0.539'
$ws.Cells.Item(8,4).ClearFormats()
$ws.Cells.Item(8,5).NumberFormat = '@'
$ws.Cells.Item(8,5).Value = '  +0.23%  '
$ws.Cells.Item(8,5).ClearFormats()

# Row 9
$ws.Cells.Item(9,4).NumberFormat = '@'
$ws.Cells.Item(9,4).Value = '2.477.39'
$ws.Cells.Item(9,4).ClearFormats()
$ws.Cells.Item(9,5).NumberFormat = '@'
$ws.Cells.Item(9,5).Value = '  +2.01%  '
$ws.Cells.Item(9,5).ClearFormats()

# Row 10
$ws.Cells.Item(10,5).NumberFormat = '@'
$ws.Cells.Item(10,5).Value = '  +0.58%  '
$ws.Cells.Item(10,5).ClearFormats()

# Row 11
$ws.Cells.Item(11,5).NumberFormat = '@'
$ws.Cells.Item(11,5).Value = '  +1.71%  '
$ws.Cells.Item(11,5).ClearFormats()

# Row 12
$ws.Cells.Item(12,5).NumberFormat = '@'
$ws.Cells.Item(12,5).Value = '  +0.79%  '
$ws.Cells.Item(12,5).ClearFormats()

# Row 13
$ws.Cells.Item(13,5).NumberFormat = '@'
$ws.Cells.Item(13,5).Value = '  +0.63%  '
$ws.Cells.Item(13,5).ClearFormats()

# Row 14
$ws.Cells.Item(14,4).NumberFormat = '@'
$ws.Cells.Item(14,4).Value = '28.64'
$ws.Cells.Item(14,4).ClearFormats()
$ws.Cells.Item(14,5).NumberFormat = '@'
$ws.Cells.Item(14,5).Value = '  +4.61%  '
$ws.Cells.Item(14,5).ClearFormats()

# Row 15
$ws.Cells.Item(15,5).NumberFormat = '@'
$ws.Cells.Item(15,5).Value = '  +1.40%  '
$ws.Cells.Item(15,5).ClearFormats()

# Row 16
$ws.Cells.Item(16,4).NumberFormat = '@'
$ws.Cells.Item(16,4).Value = '2.927.96'
$ws.Cells.Item(16,4).ClearFormats()
$ws.Cells.Item(16,5).NumberFormat = '@'
$ws.Cells.Item(16,5).Value = '  +1.53%  '
$ws.Cells.Item(16,5).ClearFormats()

# Row 17
$ws.Cells.Item(17,4).NumberFormat = '@'
$ws.Cells.Item(17,4).Value = '63.103.25'
$ws.Cells.Item(17,4).ClearFormats()
$ws.Cells.Item(17,5).NumberFormat = '@'
$ws.Cells.Item(17,5).Value = '  +1.28%  '
$ws.Cells.Item(17,5).ClearFormats()

# Row 18
$ws.Cells.Item(18,4).NumberFormat = '@'
$ws.Cells.Item(18,4).Value = '2.476.92'
$ws.Cells.Item(18,4).ClearFormats()
$ws.Cells.Item(18,5).NumberFormat = '@'
$ws.Cells.Item(18,5).Value = '  +1.67%  '
$ws.Cells.Item(18,5).ClearFormats()

# Row 19
$ws.Cells.Item(19,4).NumberFormat = '@'
$ws.Cells.Item(19,4).Value = '8.15'
$ws.Cells.Item(19,4).ClearFormats()
$ws.Cells.Item(19,5).NumberFormat = '@'
$ws.Cells.Item(19,5).Value = '  +3.52%  '
$ws.Cells.Item(19,5).ClearFormats()

# Row 20
$ws.Cells.Item(20,4).NumberFormat = '@'
$ws.Cells.Item(20,4).Value = '11.02'
$ws.Cells.Item(20,4).ClearFormats()
$ws.Cells.Item(20,5).NumberFormat = '@'
$ws.Cells.Item(20,5).Value = '  +1.17%  '
$ws.Cells.Item(20,5).ClearFormats()

# Row 21
$ws.Cells.Item(21,4).NumberFormat = '@'
$ws.Cells.Item(21,4).Value = '329.90'
$ws.Cells.Item(21,4).ClearFormats()
$ws.Cells.Item(21,5).NumberFormat = '@'
$ws.Cells.Item(21,5).Value = '  +0.95%  '
$ws.Cells.Item(21,5).ClearFormats()

# Row 22
$ws.Cells.Item(22,5).NumberFormat = '@'
$ws.Cells.Item(22,5).Value = '  +9.17%  '
$ws.Cells.Item(22,5).ClearFormats()

# Row 23
$ws.Cells.Item(23,4).NumberFormat = '@'
$ws.Cells.Item(23,4).Value = '4.13'
$ws.Cells.Item(23,4).ClearFormats()

# Row 24
$ws.Cells.Item(24,5).NumberFormat = '@'
$ws.Cells.Item(24,5).Value = '  +0.32%  '
$ws.Cells.Item(24,5).ClearFormats()

# Row 25
$ws.Cells.Item(25,2).NumberFormat = '@'
$ws.Cells.Item(25,2).Value = 'Aptos'
$ws.Cells.Item(25,2).ClearFormats()
$ws.Cells.Item(25,3).NumberFormat = '@'
$ws.Cells.Item(25,3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(25,3).ClearFormats()
$ws.Cells.Item(25,4).NumberFormat = '@'
$ws.Cells.Item(25,4).Value = '9.92'
$ws.Cells.Item(25,4).ClearFormats()
$ws.Cells.Item(25,5).NumberFormat = '@'
$ws.Cells.Item(25,5).Value = '  +16.96%  '
$ws.Cells.Item(25,5).ClearFormats()

# Row 26
$ws.Cells.Item(26,2).NumberFormat = '@'
$ws.Cells.Item(26,2).Value = 'Litecoin'
$ws.Cells.Item(26,2).ClearFormats()
$ws.Cells.Item(26,3).NumberFormat = '@'
$ws.Cells.Item(26,3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(26,3).ClearFormats()
$ws.Cells.Item(26,4).NumberFormat = '@'
$ws.Cells.Item(26,4).Value = '66.21'
$ws.Cells.Item(26,4).ClearFormats()
$ws.Cells.Item(26,5).NumberFormat = '@'
$ws.Cells.Item(26,5).Value = '  +1.18%  '
$ws.Cells.Item(26,5).ClearFormats()

# Row 27
$ws.Cells.Item(27,4).NumberFormat = '@'
$ws.Cells.Item(27,4).Value = '655.49'
$ws.Cells.Item(27,4).ClearFormats()
$ws.Cells.Item(27,5).NumberFormat = '@'
$ws.Cells.Item(27,5).Value = '  +5.69%  '
$ws.Cells.Item(27,5).ClearFormats()

# Row 28
$ws.Cells.Item(28,5).NumberFormat = '@'
$ws.Cells.Item(28,5).Value = '  +1.95%  '
$ws.Cells.Item(28,5).ClearFormats()

# Row 29
$ws.Cells.Item(29,5).NumberFormat = '@'
$ws.Cells.Item(29,5).Value = '  +2.00%  '
$ws.Cells.Item(29,5).ClearFormats()

# Row 30
$ws.Cells.Item(30,4).NumberFormat = '@'
$ws.Cells.Item(30,4).Value = '1.00'
$ws.Cells.Item(30,4).ClearFormats()
$ws.Cells.Item(30,5).NumberFormat = '@'
$ws.Cells.Item(30,5).Value = '  +264.93%  '
$ws.Cells.Item(30,5).ClearFormats()

# Row 31
$ws.Cells.Item(31,5).NumberFormat = '@'
$ws.Cells.Item(31,5).Value = '  +5.06%  '
$ws.Cells.Item(31,5).ClearFormats()

# Row 32
$ws.Cells.Item(32,4).NumberFormat = '@'
$ws.Cells.Item(32,4).Value = '8.06'
$ws.Cells.Item(32,4).ClearFormats()
$ws.Cells.Item(32,5).NumberFormat = '@'
$ws.Cells.Item(32,5).Value = '  -0.83%  '
$ws.Cells.Item(32,5).ClearFormats()

# Row 33
$ws.Cells.Item(33,5).NumberFormat = '@'
$ws.Cells.Item(33,5).Value = '  +1.33%  '
$ws.Cells.Item(33,5).ClearFormats()

# Row 34
$ws.Cells.Item(34,5).NumberFormat = '@'
$ws.Cells.Item(34,5).Value = '  -3.25%  '
$ws.Cells.Item(34,5).ClearFormats()

# Row 35
$ws.Cells.Item(35,5).NumberFormat = '@'
$ws.Cells.Item(35,5).Value = '  +4.26%  '
$ws.Cells.Item(35,5).ClearFormats()

# Row 36
$ws.Cells.Item(36,5).NumberFormat = '@'
$ws.Cells.Item(36,5).Value = '  +0.30%  '
$ws.Cells.Item(36,5).ClearFormats()

# Row 37
$ws.Cells.Item(37,5).NumberFormat = '@'
$ws.Cells.Item(37,5).Value = '  +0.84%  '
$ws.Cells.Item(37,5).ClearFormats()

# Row 38
$ws.Cells.Item(38,5).NumberFormat = '@'
$ws.Cells.Item(38,5).Value = '  +1.57%  '
$ws.Cells.Item(38,5).ClearFormats()

# Row 39
$ws.Cells.Item(39,5).NumberFormat = '@'
$ws.Cells.Item(39,5).Value = '  -0.42%  '
$ws.Cells.Item(39,5).ClearFormats()

# Row 40
$ws.Cells.Item(40,4).NumberFormat = '@'
$ws.Cells.Item(40,4).Value = '18.78'
$ws.Cells.Item(40,4).ClearFormats()
$ws.Cells.Item(40,5).NumberFormat = '@'
$ws.Cells.Item(40,5).Value = '  +1.19%  '
$ws.Cells.Item(40,5).ClearFormats()

# Row 41
$ws.Cells.Item(41,4).NumberFormat = '@'
$ws.Cells.Item(41,4).Value = '150.31'
$ws.Cells.Item(41,4).ClearFormats()
$ws.Cells.Item(41,5).NumberFormat = '@'
$ws.Cells.Item(41,5).Value = '  -1.04%  '
$ws.Cells.Item(41,5).ClearFormats()

# Row 42
$ws.Cells.Item(42,4).NumberFormat = '@'
$ws.Cells.Item(42,4).Value = '2.67'
$ws.Cells.Item(42,4).ClearFormats()
$ws.Cells.Item(42,5).NumberFormat = '@'
$ws.Cells.Item(42,5).Value = '  -2.86%  '
$ws.Cells.Item(42,5).ClearFormats()

# Row 43
$ws.Cells.Item(43,5).NumberFormat = '@'
$ws.Cells.Item(43,5).Value = '  +0.81%  '
$ws.Cells.Item(43,5).ClearFormats()

# Row 44
$ws.Cells.Item(44,5).NumberFormat = '@'
$ws.Cells.Item(44,5).Value = '  -49.05%  '
$ws.Cells.Item(44,5).ClearFormats()

# Row 45
$ws.Cells.Item(45,5).NumberFormat = '@'
$ws.Cells.Item(45,5).Value = '  +0.03%  '
$ws.Cells.Item(45,5).ClearFormats()

# Row 46
$ws.Cells.Item(46,4).NumberFormat = '@'
$ws.Cells.Item(46,4).Value = '154.98'
$ws.Cells.Item(46,4).ClearFormats()
$ws.Cells.Item(46,5).NumberFormat = '@'
$ws.Cells.Item(46,5).Value = '  +7.79%  '
$ws.Cells.Item(46,5).ClearFormats()

# Row 47
$ws.Cells.Item(47,5).NumberFormat = '@'
$ws.Cells.Item(47,5).Value = '  +3.65%  '
$ws.Cells.Item(47,5).ClearFormats()

# Row 48
$ws.Cells.Item(48,4).NumberFormat = '@'
$ws.Cells.Item(48,4).Value = '3.61'
$ws.Cells.Item(48,4).ClearFormats()
$ws.Cells.Item(48,5).NumberFormat = '@'
$ws.Cells.Item(48,5).Value = '  +0.86%  '
$ws.Cells.Item(48,5).ClearFormats()

# Row 49
$ws.Cells.Item(49,5).NumberFormat = '@'
$ws.Cells.Item(49,5).Value = '  +0.11%  '
$ws.Cells.Item(49,5).ClearFormats()

# Row 50
$ws.Cells.Item(50,5).NumberFormat = '@'
$ws.Cells.Item(50,5).Value = '  +2.10%  '
$ws.Cells.Item(50,5).ClearFormats()

# Row 51
$ws.Cells.Item(51,4).NumberFormat = '@'
$ws.Cells.Item(51,4).Value = '0.0516'
$ws.Cells.Item(51,4).ClearFormats()
$ws.Cells.Item(51,5).NumberFormat = '@'
$ws.Cells.Item(51,5).Value = '  +0.51%  '
$ws.Cells.Item(51,5).ClearFormats()
